$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.548.52"
$ws.Range("E2").Value = "  -0.18%  "

# Row 3
$ws.Range("D3").Value = "3.673.38"
$ws.Range("E3").Value = "  -1.07%  "

# Row 4
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").Value = "'616.74"
$ws.Range("E5").Value = "  -8.33%  "

# Row 6
$ws.Range("D6").Value = "'159.42"
$ws.Range("E6").Value = "  -1.54%  "

# Row 7
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("E8").Value = "  -0.58%  "

# Row 9
$ws.Range("E9").Value = "  -1.93%  "

# Row 10
$ws.Range("D10").Value = "'7.18"
$ws.Range("E10").Value = "  +1.28%  "

# Row 11
$ws.Range("E11").Value = "  -1.27%  "

# Row 12
$ws.Range("D12").Value = "'0.0000229"
$ws.Range("E12").Value = "  -2.95%  "

# Row 13
$ws.Range("D13").Value = "4.289.03"
$ws.Range("E13").Value = "  -1.23%  "

# Row 14
$ws.Range("D14").Value = "'32.43"
$ws.Range("E14").Value = "  -1.45%  "

# Row 15
$ws.Range("D15").Value = "3.671.17"
$ws.Range("E15").Value = "  -2.13%  "

# Row 16
$ws.Range("D16").Value = "69.593.28"
$ws.Range("E16").Value = "  -0.17%  "

# Row 17
$ws.Range("E17").Value = "  +0.53%  "

# Row 18
$ws.Range("D18").Value = "'6.50"
$ws.Range("E18").Value = "  -0.43%  "

# Row 19
$ws.Range("D19").Value = "'15.88"
$ws.Range("E19").Value = "  -2.80%  "

# Row 20
$ws.Range("D20").Value = "'10.30"
$ws.Range("E20").Value = "  +4.94%  "

# Row 21
$ws.Range("D21").Value = "'469.50"
$ws.Range("E21").Value = "  -1.00%  "

# Row 22
$ws.Range("D22").Value = "'0.648"
$ws.Range("E22").Value = "  -1.06%  "

# Row 23
$ws.Range("D23").Value = "'79.43"
$ws.Range("E23").Value = "  -1.34%  "

# Row 24
$ws.Range("D24").Value = "3.817.19"
$ws.Range("E24").Value = "  -1.16%  "

# Row 26
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "'0.0000122"
$ws.Range("E26").Value = "  -4.60%  "

# Row 27
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").Value = "'11.06"
$ws.Range("E27").Value = "  +0.56%  "

# Row 28
$ws.Range("D28").Value = "'8.72"
$ws.Range("E28").Value = "  -4.63%  "

# Row 29
$ws.Range("D29").Value = "'2.59"
$ws.Range("E29").Value = "  -3.86%  "

# Row 30
$ws.Range("D30").Value = "'1.67"
$ws.Range("E30").Value = "  -4.11%  "

# Row 31
$ws.Range("E31").Value = "  -0.13%  "

# Row 32
$ws.Range("E32").Value = "  -2.16%  "

# Row 33
$ws.Range("D33").Value = "'26.62"
$ws.Range("E33").Value = "  -1.27%  "

# Row 34
$ws.Range("E34").Value = "  -2.86%  "

# Row 35
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "'6.38"
$ws.Range("E35").Value = "  -3.42%  "

# Row 36
$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.670.48"
$ws.Range("E36").Value = "  -0.90%  "

# Row 37
$ws.Range("D37").Value = "'8.29"
$ws.Range("E37").Value = "  -3.38%  "

# Row 39
$ws.Range("D39").Value = "'178.10"
$ws.Range("E39").Value = "  +2.47%  "

# Row 40
$ws.Range("D40").Value = "'0.999"
$ws.Range("E40").Value = "  -0.19%  "

# Row 41
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "'5.78"
$ws.Range("E41").Value = "  -5.35%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'2.21"
$ws.Range("E42").Value = "  -1.69%  "

# Row 43
$ws.Range("D43").Value = "'0.0890"
$ws.Range("E43").Value = "  -2.61%  "

# Row 44
$ws.Range("D44").Value = "'0.926"
$ws.Range("E44").Value = "  -1.69%  "

# Row 45
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").Value = "'29.22"
$ws.Range("E45").Value = "  +5.21%  "

# Row 46
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'46.80"
$ws.Range("E46").Value = "  -0.62%  "

# Row 47
$ws.Range("D47").Value = "'2.70"
$ws.Range("E47").Value = "  -3.09%  "

# Row 48
$ws.Range("D48").Value = "'7.86"
$ws.Range("E48").Value = "  -0.41%  "

# Row 49
$ws.Range("D49").Value = "'0.000263"
$ws.Range("E49").Value = "  -7.48%  "

# Row 50
$ws.Range("D50").Value = "'1.04"
$ws.Range("E50").Value = "  -4.70%  "

# Row 51
$ws.Range("B51").Value = "ONDO"
$ws.Range("C51").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D51").Value = "'1.21"
$ws.Range("E51").Value = "  -6.47%  "
